# edit.ps1
# Applies the "form A e form B" update to the engagement analysis sheet:
#  - Bumps the expected Form 4 2025 count from 7 to 8 (header text)
#  - Updates a handful of Form 4 2025 counts (H column) that were recounted
#  - Recomputes the "Nivel de Engajamento" (I column) label - and its
#    matching color - for the rows whose level flipped between Alto/Medio
#  - Fixes the Form 2 / Form 3 status for Luana / Ceu Azul (row 56), which
#    were marked Ausente but should be Enviado, bumping that row to Alto

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors used by the "Nivel de Engajamento" / status cells.
# NOTE: Range.Interior.Color expects a BGR-ordered value (like VBA RGB()),
# so these are the BGR-packed forms of the sheet's RGB hex colors.
$colorAlto   = 0x006400   # dark green  (RGB 00 64 00 -> BGR is the same)
$colorMedio  = 0x07C1FF   # amber       (RGB FF C1 07 -> BGR 07 C1 FF)
$colorBaixo  = 0xCEC7FF   # pink        (RGB FF C7 CE -> BGR CE C7 FF)
$colorEnviado = 0xCEEFC6  # light green (RGB C6 EF CE -> BGR CE EF C6)
$colorAusente = 0xCEC7FF  # pink        (RGB FF C7 CE -> BGR CE C7 FF)

# --- Header: expected Form 4 2025 submissions moved from 7 to 8 ---
$ws.Range("H1").Value = "Form 4 2025" + [char]10 + "(Esperado: 8)"

# --- Form 4 2025 (H column) recounted values ---
$hUpdates = @{
    12 = 7
    14 = 8
    35 = 7
    40 = 7
    57 = 7
}
foreach ($row in $hUpdates.Keys) {
    $ws.Cells.Item($row, 8).Value = $hUpdates[$row]
}

# --- Engagement level (I column) label changes, Alto <-> Medio ---
$iUpdates = @{
    2  = "Médio"
    4  = "Médio"
    9  = "Médio"
    16 = "Médio"
    29 = "Médio"
    31 = "Médio"
    42 = "Médio"
    56 = "Alto"
}
foreach ($row in $iUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 9)
    $cell.Value = $iUpdates[$row]
    if ($iUpdates[$row] -eq "Alto") {
        $cell.Interior.Color = $colorAlto
    } else {
        $cell.Interior.Color = $colorMedio
    }
}

# --- Row 56 (Luana / Ceu Azul): Form 2 and Form 3 actually sent ---
$e56 = $ws.Cells.Item(56, 5)
$e56.Value = "Enviado"
$e56.Interior.Color = $colorEnviado

$f56 = $ws.Cells.Item(56, 6)
$f56.Value = "Enviado"
$f56.Interior.Color = $colorEnviado
